$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:AGRITECH"
$ws.Range("C2").Value = "NSE:GILLETTE"
$ws.Range("D2").Value = "NSE:PIIND"
$ws.Range("F2").Value = "NSE:AMBUJACEM"
$ws.Range("B3").Value = "NSE:AIAENG"
$ws.Range("C3").Value = "NSE:GRAPHITE"
$ws.Range("D3").ClearContents()
$ws.Range("F3").Value = "NSE:NAVINFLUOR"
$ws.Range("B4").Value = "NSE:BBL"
$ws.Range("C4").Value = "NSE:HEG"
$ws.Range("D4").ClearContents()
$ws.Range("F4").Value = "NSE:PIIND"
$ws.Range("B5").Value = "NSE:CHENNPETRO"
$ws.Range("C5").Value = "NSE:INFY"
$ws.Range("D5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("B6").Value = "NSE:CLEAN"
$ws.Range("C6").Value = "NSE:ITBEES"
$ws.Range("D6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("B7").Value = "NSE:CSBBANK"
$ws.Range("C7").Value = "NSE:LOYALTEX"
$ws.Range("D7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("B8").Value = "NSE:DHANUKA"
$ws.Range("C8").Value = "NSE:LYPSAGEMS"
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("B9").Value = "NSE:GOCOLORS"
$ws.Range("D9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("B10").Value = "NSE:GRAVITA"
$ws.Range("D10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("B11").Value = "NSE:HMVL"
$ws.Range("D11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("B12").Value = "NSE:HUDCO"
$ws.Range("D12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("B13").Value = "NSE:JIOFIN"
$ws.Range("D13").ClearContents()
$ws.Range("B14").Value = "NSE:JYOTHYLAB"
$ws.Range("D14").ClearContents()
$ws.Range("B15").Value = "NSE:KALYANKJIL"
$ws.Range("D15").ClearContents()
$ws.Range("B16").Value = "NSE:KEI"
$ws.Range("D16").ClearContents()
$ws.Range("B17").Value = "NSE:LOKESHMACH"
$ws.Range("D17").ClearContents()
$ws.Range("B18").Value = "NSE:MSUMI"
$ws.Range("B19").Value = "NSE:NAM-INDIA"
$ws.Range("B20").Value = "NSE:NAVINFLUOR"
$ws.Range("B21").Value = "NSE:NESCO"
$ws.Range("B22").Value = "NSE:NEWGEN"
$ws.Range("B23").Value = "NSE:NH"
$ws.Range("B24").Value = "NSE:NLCINDIA"
$ws.Range("B25").Value = "NSE:PENIND"
$ws.Range("B26").Value = "NSE:PIIND"
$ws.Range("B27").Value = "NSE:RATEGAIN"
$ws.Range("B28").Value = "NSE:RITES"

$ws.Rows("29:31").Delete()
